# Fruta / hortaliza, semanal
# Insert 5 new weekly price rows (for date serial 44504) right before the
# existing row 484, pushing the remaining data (old rows 484:567) down to
# 489:572.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing rows 484:567 down by 5 so we have room for the new rows.
$ws.Rows("484:488").Insert()

# Constant columns shared by every row in this data block.
$mercadoId   = 8
$mercado     = "Terminal La Palmera de La Serena"
$region      = "Coquimbo"
$codreg      = 4
$tipo        = "Fruta"
$productoId  = 100102
$producto    = "Cítricos"
$categoriaId = 100102005
$categoria   = "Naranja"
$unidad      = "$/bins (400 kilos)"
$origen      = "Provincia de Limarí"
$kgUnidad    = 400

# New rows: date (serial), variedad, calidad, volumen, min, max, promedio, precio/kg
$newRows = @(
    @{ Row = 484; Fecha = 44504; Variedad = "Cara cara";  Calidad = "Primera"; Volumen = 20; Min = 185000; Max = 190000; Prom = 187500; PrecioKg = 469 },
    @{ Row = 485; Fecha = 44504; Variedad = "Lane Late";  Calidad = "Primera"; Volumen = 20; Min = 155000; Max = 160000; Prom = 157500; PrecioKg = 394 },
    @{ Row = 486; Fecha = 44504; Variedad = "Lane Late";  Calidad = "Segunda"; Volumen = 20; Min = 135000; Max = 140000; Prom = 137500; PrecioKg = 344 },
    @{ Row = 487; Fecha = 44504; Variedad = "Navel Late"; Calidad = "Primera"; Volumen = 16; Min = 155000; Max = 160000; Prom = 157500; PrecioKg = 394 },
    @{ Row = 488; Fecha = 44504; Variedad = "Navel Late"; Calidad = "Segunda"; Volumen = 20; Min = 135000; Max = 140000; Prom = 137500; PrecioKg = 344 }
)

foreach ($r in $newRows) {
    $row = $r.Row
    $ws.Cells.Item($row, 1).Value2  = $mercadoId
    $ws.Cells.Item($row, 2).Value2  = $mercado
    $ws.Cells.Item($row, 3).Value2  = $region
    $ws.Cells.Item($row, 4).Value2  = $r.Fecha
    $ws.Cells.Item($row, 5).Value2  = $codreg
    $ws.Cells.Item($row, 6).Value2  = $tipo
    $ws.Cells.Item($row, 7).Value2  = $productoId
    $ws.Cells.Item($row, 8).Value2  = $producto
    $ws.Cells.Item($row, 9).Value2  = $categoriaId
    $ws.Cells.Item($row, 10).Value2 = $categoria
    $ws.Cells.Item($row, 11).Value2 = $r.Variedad
    $ws.Cells.Item($row, 12).Value2 = $r.Calidad
    $ws.Cells.Item($row, 13).Value2 = $r.Volumen
    $ws.Cells.Item($row, 14).Value2 = $r.Min
    $ws.Cells.Item($row, 15).Value2 = $r.Max
    $ws.Cells.Item($row, 16).Value2 = $r.Prom
    $ws.Cells.Item($row, 17).Value2 = $unidad
    $ws.Cells.Item($row, 18).Value2 = $origen
    $ws.Cells.Item($row, 19).Value2 = $r.PrecioKg
    $ws.Cells.Item($row, 20).Value2 = $kgUnidad
}
